$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trailing "Friend N" name+check column pairs that have no name entered were
# left over from filling the check formulas (columns E/G/I/K/M/O) all the way
# down to row 31. Clear those leftover trailing formula cells so each row
# only carries cells through its last populated "Friend" column.
$ranges = @(
    "L3:O3",
    "N4:O4",
    "N5:O5",
    "H7:O7",
    "L10:O10",
    "J11:O11",
    "J12:O12",
    "L13:O13",
    "L14:O14",
    "J15:O15",
    "L16:O16",
    "L17:O17",
    "L18:O18",
    "J19:O19",
    "J20:O20",
    "L21:O21",
    "J22:O22",
    "J23:O23",
    "H24:O24",
    "H25:O25",
    "L26:O26",
    "N27:O27",
    "L28:O28",
    "N29:O29",
    "F30:O30",
    "N31:O31"
)

foreach ($r in $ranges) {
    $ws.Range($r).ClearContents()
}

# Match the author's final selection as recorded in the saved view state.
[void]$ws.Range("J11").Select()
